# Daily attendance processing - 2025-11-14 22:22:02
# Normalizes the "Recorded By" column (G) so that the "System" token is
# moved to the end of the comma-separated list instead of the front.
#   "System, user@example.com"                 -> "user@example.com, System"
#   "System, a@b.com, system"                   -> "system, a@b.com, System"
# (other values, e.g. those that do not start with "System", are left as-is)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Transform-RecordedBy {
    param([string]$s)

    $rawParts = $s.Split(",")
    $parts = @()
    foreach ($p in $rawParts) {
        $parts += $p.Trim()
    }

    $n = $parts.Count

    if ($n -eq 2 -and $parts[0].ToLower() -eq "system") {
        # "System, X" -> "X, System"
        $result = @($parts[1], $parts[0])
    } elseif ($n -eq 3 -and $parts[0].ToLower() -eq "system" -and $parts[2].ToLower() -eq "system") {
        # "System, X, system" -> "system, X, System"  (swap first/last, keep their own casing)
        $result = @($parts[2], $parts[1], $parts[0])
    } else {
        $result = $parts
    }

    return [string]::Join(", ", $result)
}

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $orig = $cell.Value2
    if ($orig -ne $null) {
        $new = Transform-RecordedBy $orig
        $cell.Value2 = $new
    }
}
